$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ngày, tháng nhập đầy ngăn, lô kho:"
# label + its MERGEFIELD, and remove it completely (including its paragraph
# mark), same as selecting the whole line in Word and pressing Delete.
$target = $null
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*nhập đầy ngăn, lô kho*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $insertPoint = $target.Range.Start
    $target.Range.Delete()

    # Word automatically drops a "_GoBack" bookmark at the location of the
    # last edit; recreate that behaviour at the start of the paragraph that
    # now begins where the deleted paragraph used to be.
    $bmRange = $d.Range($insertPoint, $insertPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
